$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (the "MuSCs" sending-cluster rows) entirely.
$ws.Rows("4:5").Delete()

# Update the numeric values for rows 2 and 3 to reflect the new TPM data.
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4472413333333334
$ws.Range("N2").Value = 1.341724
$ws.Range("O2").Value = 0.4361236687690723
$ws.Range("P2").Value = 0.4361236687690722
$ws.Range("Q2").Value = 0.1025624261231111
$ws.Range("R2").Value = 0.9230618351080001
$ws.Range("S2").Value = 0.4361236687690723
$ws.Range("T2").Value = 0.4361236687690722

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5782506666666666
$ws.Range("N3").Value = 1.734752
$ws.Range("O3").Value = 0.5638763312309278
$ws.Range("P3").Value = 0.5638763312309278
$ws.Range("Q3").Value = 0.1326057921315555
$ws.Range("R3").Value = 1.193452129184
$ws.Range("S3").Value = 0.5638763312309278
$ws.Range("T3").Value = 0.5638763312309278
